$d = $word.ActiveDocument

# Change 1: "ithimbeni locwaningo" -> "ithimba locwaningo" (in the chatbot/email sentence)
$d.Content.Find.Execute(
    "Uma unemibuzo mayelana ne-chatbot noma uma kukhona okungacacile, sicela uthumele i-email ithimbeni locwaningo ku-",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uma unemibuzo mayelana ne-chatbot noma uma kukhona okungacacile, sicela uthumele i-email ithimba locwaningo ku-",
    2
) | Out-Null

# Change 2: "uzizwa ukhululekile" -> "uzizwe ukhululekile"
$d.Content.Find.Execute(
    "kuzosiza ukugcina ulwazi lwakho luyimfihlo futhi kuqinisekise ukuthi uzizwa ukhululekile uma uphendula imibuzo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "kuzosiza ukugcina ulwazi lwakho luyimfihlo futhi kuqinisekise ukuthi uzizwe ukhululekile uma uphendula imibuzo.",
    2
) | Out-Null

# Change 3: "kanti iMenenja yocwaningo nguZamakhanya" -> "kanye neMenenja yocwaningo uZamakhanya"
$d.Content.Find.Execute(
    "kanti iMenenja yocwaningo nguZamakhanya Makhanya (University of Cape Town).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "kanye neMenenja yocwaningo uZamakhanya Makhanya (University of Cape Town).",
    2
) | Out-Null

# Change 4: "noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli" -> "noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli"
$d.Content.Find.Execute(
    "Uma unemibuzo noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-",
    2
) | Out-Null

# Change 5: "idokhumenti engenhla" -> "incwadi engenhla"
$d.Content.Find.Execute(
    "Uma ufunde futhi waqonda idokhumenti engenhla, vuma kulemilayezo",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uma ufunde futhi waqonda incwadi engenhla, vuma kulemilayezo",
    2
) | Out-Null
